{"js": "// Replace the date line and each \"A\u00f7B=C, D\" answer cell with its updated value.\n// Every old value in this document is unique, so an exact-text search/replace\n// for each pair is unambiguous.\nconst replacements = [\n  [\"2024-05-30 Thursday\", \"2024-05-31 Friday\"],\n  [\"669\u00f78=83, 5\", \"520\u00f75=104, 0\"],\n  [\"632\u00f79=70, 2\", \"785\u00f73=261, 2\"],\n  [\"363\u00f72=181, 1\", \"197\u00f79=21, 8\"],\n  [\"331\u00f78=41, 3\", \"957\u00f73=319, 0\"],\n  [\"819\u00f75=163, 4\", \"956\u00f78=119, 4\"],\n  [\"641\u00f78=80, 1\", \"625\u00f75=125, 0\"],\n  [\"860\u00f78=107, 4\", \"760\u00f79=84, 4\"],\n  [\"649\u00f78=81, 1\", \"711\u00f78=88, 7\"],\n  [\"129\u00f79=14, 3\", \"139\u00f77=19, 6\"],\n  [\"786\u00f72=393, 0\", \"224\u00f77=32, 0\"],\n  [\"517\u00f76=86, 1\", \"669\u00f79=74, 3\"],\n  [\"149\u00f72=74, 1\", \"598\u00f74=149, 2\"],\n  [\"471\u00f72=235, 1\", \"792\u00f78=99, 0\"],\n  [\"234\u00f74=58, 2\", \"741\u00f73=247, 0\"],\n  [\"470\u00f72=235, 0\", \"228\u00f78=28, 4\"],\n  [\"917\u00f75=183, 2\", \"296\u00f73=98, 2\"],\n  [\"460\u00f73=153, 1\", \"450\u00f74=112, 2\"],\n  [\"906\u00f78=113, 2\", \"989\u00f76=164, 5\"],\n  [\"301\u00f76=50, 1\", \"516\u00f73=172, 0\"],\n  [\"190\u00f79=21, 1\", \"939\u00f73=313, 0\"],\n  [\"810\u00f78=101, 2\", \"462\u00f74=115, 2\"],\n  [\"683\u00f72=341, 1\", \"703\u00f74=175, 3\"],\n  [\"511\u00f73=170, 1\", \"118\u00f73=39, 1\"],\n  [\"731\u00f78=91, 3\", \"147\u00f78=18, 3\"],\n  [\"589\u00f74=147, 1\", \"106\u00f72=53, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"A\u00f7B=C, D\" answer cell with its updated\n# value. Every old value in this document is unique, so an exact-text\n# Find/Replace for each pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-05-30 Thursday\", \"2024-05-31 Friday\"),\n    @(\"669\u00f78=83, 5\", \"520\u00f75=104, 0\"),\n    @(\"632\u00f79=70, 2\", \"785\u00f73=261, 2\"),\n    @(\"363\u00f72=181, 1\", \"197\u00f79=21, 8\"),\n    @(\"331\u00f78=41, 3\", \"957\u00f73=319, 0\"),\n    @(\"819\u00f75=163, 4\", \"956\u00f78=119, 4\"),\n    @(\"641\u00f78=80, 1\", \"625\u00f75=125, 0\"),\n    @(\"860\u00f78=107, 4\", \"760\u00f79=84, 4\"),\n    @(\"649\u00f78=81, 1\", \"711\u00f78=88, 7\"),\n    @(\"129\u00f79=14, 3\", \"139\u00f77=19, 6\"),\n    @(\"786\u00f72=393, 0\", \"224\u00f77=32, 0\"),\n    @(\"517\u00f76=86, 1\", \"669\u00f79=74, 3\"),\n    @(\"149\u00f72=74, 1\", \"598\u00f74=149, 2\"),\n    @(\"471\u00f72=235, 1\", \"792\u00f78=99, 0\"),\n    @(\"234\u00f74=58, 2\", \"741\u00f73=247, 0\"),\n    @(\"470\u00f72=235, 0\", \"228\u00f78=28, 4\"),\n    @(\"917\u00f75=183, 2\", \"296\u00f73=98, 2\"),\n    @(\"460\u00f73=153, 1\", \"450\u00f74=112, 2\"),\n    @(\"906\u00f78=113, 2\", \"989\u00f76=164, 5\"),\n    @(\"301\u00f76=50, 1\", \"516\u00f73=172, 0\"),\n    @(\"190\u00f79=21, 1\", \"939\u00f73=313, 0\"),\n    @(\"810\u00f78=101, 2\", \"462\u00f74=115, 2\"),\n    @(\"683\u00f72=341, 1\", \"703\u00f74=175, 3\"),\n    @(\"511\u00f73=170, 1\", \"118\u00f73=39, 1\"),\n    @(\"731\u00f78=91, 3\", \"147\u00f78=18, 3\"),\n    @(\"589\u00f74=147, 1\", \"106\u00f72=53, 0\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 1\n\n    # wdFindContinue=1, wdReplaceAll=2\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
